# Doing Updates for Financials
# Insert a new "most recent quarter" column (D) into the LOW quarterly
# financials sheet, shifting the existing D:K data right to E:L, and
# populate the new column with the latest reported figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D. This shifts the previous D:K
# columns to E:L (Excel's native "insert column" behaviour), including
# all existing values/formatting.
$ws.Columns("D").Insert(-4160, 0)

# The freshly inserted column D starts out with the default/no style.
# Copy the number formats from the (now shifted) column E so the new
# column matches the surrounding data (date format in header rows,
# the numeric format elsewhere).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# Populate the new column D with the newly reported quarter's values.

# --- Income Statement (rows 7-35) ---
$ws.Range("D7").Value2 = 43406
$ws.Range("D8").Value2 = 17415000
$ws.Range("D9").Value2 = 11755000
$ws.Range("D10").Value2 = 5660000
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("D15").Value2 = 433000
$ws.Range("D17").Value2 = 16458000
$ws.Range("D18").Value2 = 957000
$ws.Range("D20").Value2 = 7000
$ws.Range("D21").Value2 = 1419000
$ws.Range("D22").Value2 = 160000
$ws.Range("D23").Value2 = 804000
$ws.Range("D24").Value2 = 175000
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 629000
$ws.Range("D27").Value2 = 628000
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -7000
$ws.Range("D33").Value2 = 628000
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 628000

# --- Balance Sheet (rows 38-77) ---
$ws.Range("D38").Value2 = 43406
$ws.Range("D41").Value2 = 1668000
$ws.Range("D42").Value2 = 27000
$ws.Range("D43").Value2 = 0
$ws.Range("D44").Value2 = 12365000
$ws.Range("D45").Value2 = 1078000
$ws.Range("D46").Value2 = 15138000
$ws.Range("D47").Value2 = 290000
$ws.Range("D48").Value2 = 18923000
$ws.Range("D49").Value2 = 1272000
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 1090000
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 36713000
$ws.Range("D57").Value2 = 9283000
$ws.Range("D58").Value2 = 1117000
$ws.Range("D59").Value2 = 4669000
$ws.Range("D60").Value2 = 15069000
$ws.Range("D61").Value2 = 14460000
$ws.Range("D62").Value2 = 1790000
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 31319000
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 5156000
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 5394000
$ws.Range("D77").Value2 = 0

# --- Cash Flow (rows 80-102) ---
$ws.Range("D80").Value2 = 43406
$ws.Range("D81").Value2 = 628000
$ws.Range("D83").Value2 = 455000
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 1011000
$ws.Range("D91").Value2 = -303000
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -308000
$ws.Range("D96").Value2 = -390000
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -1284000
$ws.Range("D101").Value2 = -2000
$ws.Range("D102").Value2 = -583000
